$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2 through 40 all hold the date serial 45205
# (2023-10-06) and must be bumped to 45206 (2023-10-07).
for ($row = 2; $row -le 40; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
